$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values in column D look like plain numbers (e.g. "545.19"), but
# column D stores every price as text (others contain multiple "." separators that
# are not valid numbers, e.g. "60.102.57"). Force those cells to Text format first
# so Excel does not auto-convert the numeric-looking ones to actual numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.102.57"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("D3").Value = "2.347.04"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "545.19"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").Value = "132.07"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "2.345.69"
$ws.Range("E9").Value = "  +3.19%  "
$ws.Range("E10").Value = "  +2.32%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "23.83"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "2.762.90"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").Value = "60.077.94"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "2.354.04"
$ws.Range("E18").Value = "  +4.35%  "
$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  +6.84%  "
$ws.Range("D22").Value = "313.88"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "63.39"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "7.88"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "1.35"
$ws.Range("E28").Value = "  +6.59%  "
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("D30").Value = "171.58"
$ws.Range("E31").Value = "  +11.66%  "
$ws.Range("D32").Value = "0.0₃0726"
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("E34").Value = "  +14.95%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "18.02"
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "4.16"
$ws.Range("E39").Value = "  +7.52%  "
$ws.Range("D40").Value = "321.70"
$ws.Range("E40").Value = "  +13.05%  "
$ws.Range("D41").Value = "38.13"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E42").Value = "  +3.77%  "
$ws.Range("D43").Value = "141.55"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "3.45"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "19.40"
$ws.Range("E46").Value = "  +7.99%  "
$ws.Range("D47").Value = "0.0497"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").Value = "0.562"
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "0.0₆0214"
$ws.Range("E50").Value = "  +19.74%  "
$ws.Range("D51").Value = "11.01"
$ws.Range("E51").Value = "  +0.82%  "
